$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the columns that should no longer contain data (C1:O1)
$ws.Range("C1:O1").Clear()

# Set new values for the remaining cells
$ws.Range("A1").Value = 20
$ws.Range("B1").Value = 21
